$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 123 (pushes existing rows 123:173 down to 124:174,
# carrying their data/format with them - matches the diff where every row
# from 123 on is shifted down by one and a new record appears at the top).
$ws.Rows.Item(123).Insert()

# Populate the newly inserted row 123 with the new "Acelga" price record.
$ws.Cells.Item(123, 1).Value = 5
$ws.Cells.Item(123, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(123, 3).Value = "Maule"
$ws.Cells.Item(123, 4).Value = 44466
$ws.Cells.Item(123, 5).Value = 7
$ws.Cells.Item(123, 6).Value = 100112009
$ws.Cells.Item(123, 7).Value = "Acelga"
$ws.Cells.Item(123, 8).Value = "Sin especificar"
$ws.Cells.Item(123, 9).Value = "Primera"
$ws.Cells.Item(123, 10).Value = 500
$ws.Cells.Item(123, 11).Value = 2000
$ws.Cells.Item(123, 12).Value = 2000
$ws.Cells.Item(123, 13).Value = 2000
$ws.Cells.Item(123, 14).Value = "`$/docena de atados (4 kilos)"
$ws.Cells.Item(123, 15).Value = "Región del Maule"
$ws.Cells.Item(123, 16).Value = 500
$ws.Cells.Item(123, 17).Value = 4
$ws.Cells.Item(123, 18).Value = "Hortaliza"
